$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "added this line in main branch "
$ws.Range("A3").Select()
